# Insert a new weekly data row right above the current row 15 (the first
# record for Sin especificar / 04-03-2022), shifting all subsequent rows
# down by one. This also pushes the former last data row (old row 74) into
# a brand-new row 75.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A15").EntireRow.Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Range("A15").Value = 7
$ws.Range("B15").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C15").Value = "Ñuble"
$ws.Range("D15").Value = 44624
$ws.Range("E15").Value = 16
$ws.Range("F15").Value = 100112031
$ws.Range("G15").Value = "Poroto verde"
$ws.Range("H15").Value = "Sin especificar"
$ws.Range("I15").Value = "Primera"
$ws.Range("J15").Value = 120
$ws.Range("K15").Value = 28000
$ws.Range("L15").Value = 29000
$ws.Range("M15").Value = 28500
$ws.Range("N15").Value = "`$/saco 25 kilos"
$ws.Range("O15").Value = "Región del Maule"
$ws.Range("P15").Value = 1140
$ws.Range("Q15").Value = 25
$ws.Range("R15").Value = "Hortaliza"
